$d = $word.ActiveDocument
$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ------------------------------------------------------------------
# 1) Insert the new "Steps" and "Checklist" sections (plus their
#    surrounding blank paragraphs) before the existing "Ideas"
#    heading paragraph, which currently sits at the very start of
#    the body.
# ------------------------------------------------------------------
$newSectionsXml = @"
<w:p $wns/>
<w:p $wns>
  <w:pPr>
    <w:pStyle w:val="Heading1"/>
  </w:pPr>
  <w:r>
    <w:t>Steps</w:t>
  </w:r>
</w:p>
<w:p $wns>
  <w:r>
    <w:t xml:space="preserve">All </w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve">nav buttons should link back </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:proofErr w:type="gramStart"/>
  <w:r>
    <w:t>tp</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:proofErr w:type="gramEnd"/>
</w:p>
<w:p $wns/>
<w:p $wns>
  <w:pPr>
    <w:pStyle w:val="Heading1"/>
  </w:pPr>
  <w:r>
    <w:t>Checklist</w:t>
  </w:r>
</w:p>
<w:p $wns>
  <w:r>
    <w:t xml:space="preserve">Make sure all nav headings link back to the various page </w:t>
  </w:r>
  <w:proofErr w:type="gramStart"/>
  <w:r>
    <w:t>sections</w:t>
  </w:r>
  <w:proofErr w:type="gramEnd"/>
</w:p>
<w:p $wns>
  <w:r>
    <w:t xml:space="preserve">Add a </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:proofErr w:type="gramStart"/>
  <w:r>
    <w:t>flavicon</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:proofErr w:type="gramEnd"/>
</w:p>
<w:p $wns/>
"@

$insertionPoint = $d.Range(0, 0)
$insertionPoint.InsertXML($newSectionsXml)

# ------------------------------------------------------------------
# 2) The paragraph that used to be blank right after "Ideas" is no
#    longer needed -- in the new layout "Ideas" is immediately
#    followed by "Make the skills pictures circle".  Find the
#    "Ideas" heading (now relocated after the inserted sections) and
#    drop the empty paragraph that still trails it.
# ------------------------------------------------------------------
$ideasRange = $d.Content
$ideasRange.Find.Execute("Ideas", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$ideasParagraph = $ideasRange.Paragraphs.First
$trailingBlank = $ideasParagraph.Next
$trailingBlank.Range.Delete()

# ------------------------------------------------------------------
# 3) Append the new "Keep portfolio pics square" bullet (plus a
#    trailing blank paragraph) after the existing
#    "Make the skills pictures circle" paragraph.
# ------------------------------------------------------------------
$lastParagraph = $d.Paragraphs.Last
$lastParagraph.Range.InsertParagraphAfter()
$appendTarget = $d.Paragraphs.Last

$tailXml = @"
<w:p $wns>
  <w:r>
    <w:t xml:space="preserve">Keep portfolio pics </w:t>
  </w:r>
  <w:proofErr w:type="gramStart"/>
  <w:r>
    <w:t>square</w:t>
  </w:r>
  <w:proofErr w:type="gramEnd"/>
</w:p>
<w:p $wns/>
"@

$appendTarget.Range.InsertXML($tailXml)
